$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.904.69"
$ws.Range("E2").Value = "  +1.43%  "
$ws.Range("D3").Value = "1.640.51"
$ws.Range("E3").Value = "  +1.19%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.73"
$ws.Range("E5").Value = "  +1.09%  "
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.75"
$ws.Range("E8").Value = "  +2.90%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("E10").Value = "  +0.81%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0874"
$ws.Range("E11").Value = "  -0.68%  "
$ws.Range("D12").Value = "1.873.79"
$ws.Range("E12").Value = "  +1.20%  "
$ws.Range("D13").Value = "1.640.95"
$ws.Range("E13").Value = "  +0.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.576"
$ws.Range("E14").Value = "  +4.79%  "
$ws.Range("E15").Value = "  +1.33%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.16"
$ws.Range("E16").Value = "  +1.30%  "
$ws.Range("D17").Value = "27.906.95"
$ws.Range("E17").Value = "  +1.51%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "231.93"
$ws.Range("E18").Value = "  +0.79%  "
$ws.Range("D19").Value = "0.0₃0725"
$ws.Range("E19").Value = "  +1.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.61"
$ws.Range("E20").Value = "  +0.91%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.51"
$ws.Range("E21").Value = "  +10.85%  "
$ws.Range("E23").Value = "  +1.24%  "
$ws.Range("E24").Value = "  -2.53%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.92"
$ws.Range("E25").Value = "  +1.42%  "
$ws.Range("E26").Value = "  +1.19%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.75"
$ws.Range("E27").Value = "  +1.43%  "
$ws.Range("B28").Value = "Stellar"
$ws.Range("C28").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.112"
$ws.Range("E28").Value = "  +0.75%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("E30").Value = "  +1.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0485"
$ws.Range("E31").Value = "  +0.34%  "
$ws.Range("E32").Value = "  +2.22%  "
$ws.Range("D33").Value = "1.424.06"
$ws.Range("E33").Value = "  -2.81%  "
$ws.Range("E34").Value = "  +2.20%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.57"
$ws.Range("E35").Value = "  +1.80%  "
$ws.Range("E36").Value = "  +0.44%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.893"
$ws.Range("E37").Value = "  +2.69%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0167"
$ws.Range("E38").Value = "  +0.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.923"
$ws.Range("E39").Value = "  -2.32%  "
$ws.Range("E40").Value = "  +1.10%  "
$ws.Range("E41").Value = "  +2.26%  "
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "67.28"
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("E44").Value = "  -1.40%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.82"
$ws.Range("E45").Value = "  +4.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.46"
$ws.Range("E46").Value = "  +2.69%  "
$ws.Range("E47").Value = "  +0.34%  "
$ws.Range("D48").Value = "1.782.55"
$ws.Range("E49").Value = "  +1.55%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.101"
$ws.Range("E50").Value = "  +1.12%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0506"
$ws.Range("E51").Value = "  +0.70%  "
